$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (the row with "iaest-measure:..." values),
# shifting existing rows 2-4 down to rows 3-5.
$ws.Rows.Item(2).Insert()

# Fill the newly inserted row 2 with the slug-style identifiers that relate
# to the header row (A1:E1), used to build SKOS hierarchies between columns.
$ws.Range("A2").Value = "horas-trabajadas"
$ws.Range("B2").Value = "personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "sector-actividad-descripcion"
$ws.Range("D2").Value = "aragon"
$ws.Range("E2").Value = "sector-actividad-codigo"
